# Fruta / hortaliza, semanal
# A new weekly price-report entry was inserted above the existing row 203,
# pushing the old rows 203:235 down to 204:236 (dimension grows from
# A1:T235 to A1:T236). Fill the newly inserted row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 203; Excel shifts rows 203:235 down to 204:236,
# carrying their formatting (including the D-column date style) with them.
$ws.Rows(203).Insert()

$ws.Range("A203").Value = 5
$ws.Range("B203").Value = "Macroferia Regional de Talca"
$ws.Range("C203").Value = "Maule"
$ws.Range("D203").Value = 44522
$ws.Range("E203").Value = 7
$ws.Range("F203").Value = "Fruta"
$ws.Range("G203").Value = 100102
$ws.Range("H203").Value = "Cítricos"
$ws.Range("I203").Value = 100102004
$ws.Range("J203").Value = "Mandarina"
$ws.Range("K203").Value = "Murcott"
$ws.Range("L203").Value = "Primera"
$ws.Range("M203").Value = 150
$ws.Range("N203").Value = 6000
$ws.Range("O203").Value = 6000
$ws.Range("P203").Value = 6000
$ws.Range("Q203").Value = "$/bandeja 10 kilos"
$ws.Range("R203").Value = "Provincia de Limarí"
$ws.Range("S203").Value = 600
$ws.Range("T203").Value = 10
